# Auto: Update ETF Data
# Reset "Share Change" (G) and "Net Amount" (H) to 0 for the rows whose
# values changed in the source data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6, 10, 13, 14, 17, 25, 28, 33, 34, 39, 40, 49, 52)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
}
